$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" everywhere it appears ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the "zh-cn"/"de-de" status columns (Overview!E:F and the Status column on each locale sheet) ---
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101848602295
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101848602295

$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101848602295

$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101848602295
